$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record as row 13; this pushes the existing
# rows 13..58 down to 14..59 (dimension grows from R58 to R59).
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new record's data. The columns that
# are constant for every record in this sheet (A, B, C, E, F, G, H, I, Q, R)
# are set explicitly too, since Insert() leaves the new row blank.
$ws.Range("A13").Value2 = 3
$ws.Range("B13").Value = "Femacal de La Calera"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value2 = 44676
$ws.Range("E13").Value2 = 5
$ws.Range("F13").Value2 = 100112022
$ws.Range("G13").Value = "Arveja Verde"
$ws.Range("H13").Value = "Perfection"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value2 = 73
$ws.Range("K13").Value2 = 23000
$ws.Range("L13").Value2 = 24000
$ws.Range("M13").Value2 = 23479
$ws.Range("N13").Value = "$/malla 25 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value2 = 939
$ws.Range("Q13").Value2 = 25
$ws.Range("R13").Value = "Hortaliza"
